$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Step1 - Input Data")

# Update the Personal Required Rate of Return from 6.86% to 6.87%.
# The cell already carries a percentage number format (0%) but its
# stored value is a literal text string "6.86%", not a real number.
# Assigning the new text directly would make Excel auto-parse it into
# a numeric percentage (0.0687) because of the cell's number format,
# so we temporarily switch the cell to Text format while writing the
# new value, then restore the original percentage number format
# afterwards. This keeps the cell's stored type/value as text "6.87%",
# matching the original authoring pattern.
$cell = $ws.Range("B4")
$originalFormat = $cell.NumberFormat
$cell.NumberFormat = "@"
$cell.Value2 = "6.87%"
$cell.NumberFormat = $originalFormat

# Recalculate dependent formulas (Step2 - Projection sheet references
# this cell and several downstream cells derive from it).
$excel.CalculateFullRebuild()
